# ProjectsTracker.xlsx update:
#  - insert a "skills" column (new col B) between title and description
#  - insert a "location" column (new col D) between description and link
#  - fill in the skills/location values for the two existing project rows
#  - nudge the remembered selection from E21 -> E18 (matches the new, shorter layout)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- insert the two new columns, shifting description/link out of the way ---
# Before: A=title B=description C=link
$ws.Columns.Item(2).EntireColumn.Insert()   # new, empty B -> description shifts B->C, link C->D
$ws.Columns.Item(4).EntireColumn.Insert()   # new, empty D -> link shifts D->E
# After:  A=title B=(new) C=description D=(new) E=link

# --- fill the "skills" column (B) top to bottom ---
$ws.Range("B1").Value = "skills"
$ws.Range("B2").Value = "Computer Vision, OpenCV, Java, Image Processing"
$ws.Range("B3").Value = "Javascript, HTML, CSS, LLM, API"

# --- fill the "location" column (D) top to bottom ---
$ws.Range("D1").Value = "location"
$ws.Range("D2").Value = "View on Github"
$ws.Range("D3").Value = "View on Github"

# --- column widths for the two new columns (best-fit-ish sizing) ---
$ws.Columns.Item(2).ColumnWidth = 41.166666
$ws.Columns.Item(4).ColumnWidth = 12

# --- restore the selected cell Excel remembered on save ---
$ws.Range("E18").Select() | Out-Null
